$wb = $excel.ActiveWorkbook

# --- Sheet1: update Base Image / Thumbnail Image values for the test rows ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2 (SKU01 / test)
$ws1.Range("F2").Value = "test.jpg"
$ws1.Range("G2").Value = "test_thumb.jpg"

# Row 3 (SKU02 / test2)
$ws1.Range("F3").Value = "test2.jpg"
$ws1.Range("G3").Value = "test2_thumb.jpg"

# Row 4 (fdsaf / dasfd) - thumbnail updated to the new naming convention
$ws1.Range("G4").Value = "fdsafd_thumb.jpg"

# Select a neutral cell on Sheet1 (matches the post-edit workbook: Sheet1 is no
# longer the tab shown with a live selection highlight baked into its view)
$ws1.Range("I2").Select() | Out-Null

# --- Sheet2 becomes the active/selected sheet ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate() | Out-Null
